# Apply updated cryptocurrency price/volume data to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.924.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.73%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.810.09"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.96%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.01%  "

$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4624"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.86%  "

$ws.Range("E8").Value = "  -1.79%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07376"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8771"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.32%  "

$ws.Range("E11").Value = "  -1.77%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.767.65"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.35%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.360"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.33%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.516"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.15%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07041"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.20%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008713"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.97%  "

$ws.Range("E19").Value = "  +0.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.74"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.941.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.323"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.76%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.93%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.075.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.51%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.899"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.52%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.26%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.23%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.160"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.48%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.350"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.12%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08900"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7560"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.65%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.456"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.45%  "

$ws.Range("E35").Value = "  -0.65%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9997"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.00%  "

$ws.Range("E37").Value = "  -0.71%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01972"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.41%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05254"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.418"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.28%  "

$ws.Range("E41").Value = "  +1.82%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5335"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.231"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1664"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.519"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.85%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4989"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.70%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.674"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.84%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9996"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06293"
$ws.Range("D51").Style = "Normal"
